$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

$ws.Range("A20").Value = "26-07-2018"
$ws.Range("B20").Value = "26-07-2018"
$ws.Range("C20").Value = "parametric_ao_NA"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "PLA"
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 0.2
$ws.Range("I20").Value = "NA"

$ws.Range("A20:I20").HorizontalAlignment = -4108

$ws.Range("C23").Select()
